$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 47, pushing the existing rows 47-50 down to 48-51.
$ws.Rows(47).Insert()

# Populate the newly inserted row 47 with the waiting-for-players localization entry.
$ws.Cells.Item(47, 1).Value = "UI_GAME_WAITING_FOR_PLAYERS"
$ws.Cells.Item(47, 2).Value = "Waiting for other players…"
$ws.Cells.Item(47, 3).Value = "XXXX"
$ws.Cells.Item(47, 4).Value = "XXXX"
$ws.Cells.Item(47, 5).Value = "XXXX"

# Match the author's final selection state.
[void]$ws.Range("B48").Select()
